# Reorganize evaluation sheet columns:
# Add a new "unclear" column (E) and move rows where both C ("correct")
# and D ("wrong") were marked with "x" into this new column instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New header for column E
$ws.Range("E3").Value = "unclear"

# Rows that had both C and D marked "x" (ambiguous) -> move to E, clear C & D
$rows = @(28, 30, 33, 35, 42)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "x"   # column E
    $ws.Cells.Item($r, 3).ClearContents() # column C
    $ws.Cells.Item($r, 4).ClearContents() # column D
}

$ws.Range("F27").Select()
